$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column (H): copy the header formatting from the neighboring
# header cell (G1) so H1 reuses the existing bold/centered/bordered style,
# then set its text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Save flag values for H2:H7
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 0
